# Add a new "Monster Hunter 3 G" column (column V) to the monster list sheet.
# Monster Hunter 3 G shares its monster roster with Monster Hunter 3 Ultimate
# (column I), so the new column is populated with the same monster names.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new column.
$ws.Cells.Item(1, 22).Value2 = "Monster Hunter 3 G"

# Copy the Monster Hunter 3 Ultimate monster list (column I, rows 2-74)
# into the new Monster Hunter 3 G column (column V).
for ($r = 2; $r -le 74; $r++) {
    $ws.Cells.Item($r, 22).Value2 = $ws.Cells.Item($r, 9).Value2
}

# Give the new column a width matching the other monster-name columns.
$ws.Columns.Item(22).ColumnWidth = 26

# Move the selection/view roughly where the author left it after the edit.
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 10
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("K2").Select()
